$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 30385
$ws.Range("J43").Value = 18573.166
$ws.Range("L43").Value = 18573.166
$ws.Range("N43").Value = -18711.166
$ws.Range("H127").Value = 333504
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 333504
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 1000512
$ws.Range("M127").Value = $null
$ws.Range("N127").Value = -1010432
$ws.Range("H134").Value = 119999.5
$ws.Range("J134").Value = 119999.5
$ws.Range("L134").Value = 119999.5
$ws.Range("N134").Value = -130139.5
$ws.Range("H135").Value = 1180.2142
$ws.Range("I135").Value = 963.3077
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 8669.7693
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -6134.7693
$ws.Range("N135").Value = -41070
$ws.Range("H137").Value = 2302.889
$ws.Range("J137").Value = 1937.5714
$ws.Range("L137").Value = 5812.7142
$ws.Range("N137").Value = -10912.7142
$ws.Range("H138").Value = 3543.7632
$ws.Range("I138").Value = 1425.5625
$ws.Range("K138").Value = 4276.6875
$ws.Range("M138").Value = 863.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 28251
$ws.Range("I10").Value = 4999
$ws.Range("J10").Value = 36001.668
$ws.Range("K10").Value = 4999
$ws.Range("L10").Value = 36001.668
$ws.Range("M10").Value = -4829
$ws.Range("N10").Value = -36341.668
$ws.Range("H23").Value = 91666.664
$ws.Range("J23").Value = 125000
$ws.Range("L23").Value = 125000
$ws.Range("N23").Value = -125518
$ws.Range("H32").Value = 4026.1707
$ws.Range("I32").Value = 3144.6365
$ws.Range("K32").Value = 3144.6365
$ws.Range("M32").Value = -2857.6365
$ws.Range("H97").Value = 2716
$ws.Range("I97").Value = 2174.7368
$ws.Range("K97").Value = 2174.7368
$ws.Range("M97").Value = -1678.7368
$ws.Range("H110").Value = 10536.782
$ws.Range("I110").Value = 16316.333
$ws.Range("K110").Value = 16316.333
$ws.Range("M110").Value = -14271.333
$ws.Range("H132").Value = 5106.1113
$ws.Range("I132").Value = 3492.5
$ws.Range("J132").Value = 8333.333000000001
$ws.Range("K132").Value = 10477.5
$ws.Range("L132").Value = 24999.999
$ws.Range("M132").Value = -7947.5
$ws.Range("N132").Value = -30059.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 40000
$ws.Range("I26").Value = 40000
$ws.Range("K26").Value = 40000
$ws.Range("M26").Value = -39708
$ws.Range("H86").Value = 5242.357
$ws.Range("I86").Value = 3309.6
$ws.Range("J86").Value = 7472.4614
$ws.Range("K86").Value = 3309.6
$ws.Range("L86").Value = 7472.4614
$ws.Range("M86").Value = -2186.6
$ws.Range("N86").Value = -9718.4614
$ws.Range("H89").Value = 5242.357
$ws.Range("I89").Value = 3309.6
$ws.Range("J89").Value = 7472.4614
$ws.Range("K89").Value = 16548
$ws.Range("L89").Value = 37362.307
$ws.Range("M89").Value = -10932
$ws.Range("N89").Value = -48594.307
$ws.Range("H126").Value = 70001
$ws.Range("J126").Value = 70001
$ws.Range("L126").Value = 70001
$ws.Range("N126").Value = -79881

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47495.22
$ws.Range("J31").Value = 11752.167
$ws.Range("L31").Value = 11752.167
$ws.Range("N31").Value = -12342.167
$ws.Range("H34").Value = 47495.22
$ws.Range("J34").Value = 11752.167
$ws.Range("L34").Value = 11752.167
$ws.Range("N34").Value = -12156.167
$ws.Range("H58").Value = 2810
$ws.Range("I58").Value = 2911.25
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 2911.25
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -2708.25
$ws.Range("N58").Value = -2406
$ws.Range("H99").Value = 2627.3125
$ws.Range("I99").Value = 2364.5454
$ws.Range("J99").Value = 3205.4
$ws.Range("K99").Value = 2364.5454
$ws.Range("L99").Value = 3205.4
$ws.Range("M99").Value = -866.5454
$ws.Range("N99").Value = -6201.4
$ws.Range("H126").Value = 2627.3125
$ws.Range("I126").Value = 2364.5454
$ws.Range("J126").Value = 3205.4
$ws.Range("K126").Value = 7093.6362
$ws.Range("L126").Value = 9616.200000000001
$ws.Range("M126").Value = -4623.6362
$ws.Range("N126").Value = -14556.2
$ws.Range("H132").Value = 3222.7454
$ws.Range("I132").Value = 3165.6
$ws.Range("K132").Value = 9496.799999999999
$ws.Range("M132").Value = -6966.799999999999
$ws.Range("H134").Value = 6839.1406
$ws.Range("I134").Value = 4935.1753
$ws.Range("K134").Value = 14805.5259
$ws.Range("M134").Value = -12270.5259
$ws.Range("H136").Value = 2810
$ws.Range("I136").Value = 2911.25
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 8733.75
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -6183.75
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 45.51852
$ws.Range("I2").Value = 52.444443
$ws.Range("K2").Value = 314.666658
$ws.Range("M2").Value = -201.666658
$ws.Range("H23").Value = 254.90475
$ws.Range("I23").Value = 162.1
$ws.Range("J23").Value = 339.27274
$ws.Range("K23").Value = 486.3
$ws.Range("L23").Value = 1017.81822
$ws.Range("M23").Value = -251.3
$ws.Range("N23").Value = -1487.81822
$ws.Range("H34").Value = 1275.0834
$ws.Range("J34").Value = 2141.1428
$ws.Range("L34").Value = 6423.428400000001
$ws.Range("N34").Value = -6591.428400000001
$ws.Range("H113").Value = 538.3333
$ws.Range("J113").Value = 565.7143
$ws.Range("L113").Value = 1697.1429
$ws.Range("N113").Value = -6037.1429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11411
$ws.Range("J70").Value = 13013.75
$ws.Range("L70").Value = 13013.75
$ws.Range("N70").Value = -13553.75
$ws.Range("H73").Value = 11411
$ws.Range("J73").Value = 13013.75
$ws.Range("L73").Value = 13013.75
$ws.Range("N73").Value = -14885.75
$ws.Range("H122").Value = 3757.5264
$ws.Range("I122").Value = 3574.7
$ws.Range("J122").Value = 3960.6667
$ws.Range("K122").Value = 10724.1
$ws.Range("L122").Value = 11882.0001
$ws.Range("M122").Value = -8274.099999999999
$ws.Range("N122").Value = -16782.0001
$ws.Range("H136").Value = 41250
$ws.Range("J136").Value = 41250
$ws.Range("L136").Value = 123750
$ws.Range("N136").Value = -128850

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2008.3334
$ws.Range("J22").Value = 2112.5
$ws.Range("L22").Value = 2112.5
$ws.Range("N22").Value = -2702.5
$ws.Range("H27").Value = 2008.3334
$ws.Range("J27").Value = 2112.5
$ws.Range("L27").Value = 2112.5
$ws.Range("N27").Value = -2326.5
$ws.Range("H40").Value = 4376.7896
$ws.Range("I40").Value = 3609.9333
$ws.Range("K40").Value = 3609.9333
$ws.Range("M40").Value = -3473.9333
$ws.Range("H46").Value = 1677.3636
$ws.Range("I46").Value = 1494.3334
$ws.Range("J46").Value = 2501
$ws.Range("K46").Value = 1494.3334
$ws.Range("L46").Value = 2501
$ws.Range("M46").Value = -1306.3334
$ws.Range("N46").Value = -2877
$ws.Range("H55").Value = 538.7143
$ws.Range("I55").Value = 399.75
$ws.Range("J55").Value = 724
$ws.Range("K55").Value = 399.75
$ws.Range("L55").Value = 724
$ws.Range("M55").Value = -226.75
$ws.Range("N55").Value = -1070
$ws.Range("H61").Value = 4554.8125
$ws.Range("I61").Value = 4243.8184
$ws.Range("K61").Value = 4243.8184
$ws.Range("M61").Value = -4041.8184
$ws.Range("H113").Value = 4554.8125
$ws.Range("I113").Value = 4243.8184
$ws.Range("K113").Value = 4243.8184
$ws.Range("M113").Value = -2073.8184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 32948.5
$ws.Range("I2").Value = 28975.555
$ws.Range("J2").Value = 40099.8
$ws.Range("K2").Value = 28975.555
$ws.Range("L2").Value = 40099.8
$ws.Range("M2").Value = -28863.555
$ws.Range("N2").Value = -40323.8
$ws.Range("H64").Value = 12120
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 12120
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H126").Value = 3916.9473
$ws.Range("I126").Value = 2524.3845
$ws.Range("J126").Value = 6934.1665
$ws.Range("K126").Value = 7573.1535
$ws.Range("L126").Value = 20802.4995
$ws.Range("M126").Value = -5103.1535
$ws.Range("N126").Value = -25742.4995
$ws.Range("H132").Value = 5649.05
$ws.Range("I132").Value = 5577.9473
$ws.Range("K132").Value = 16733.8419
$ws.Range("M132").Value = -14203.8419
